# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" and "全部类型" sheets, which carry the same source data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 367
    3 = 1227
    4 = 1439
    5 = 55
    6 = 6109
    7 = 102
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
